$wb = $excel.ActiveWorkbook

# Rename sheets (order matches rId1..rId5 -> sheet1..sheet5)
$wb.Worksheets.Item(1).Name = "GNG_TO-16502912981973264"
$wb.Worksheets.Item(2).Name = "NB_TO-16502912999326546"
$wb.Worksheets.Item(3).Name = "RS_TO-1650291299934654"
$wb.Worksheets.Item(4).Name = "TOL_TO-16502912999959936"
$wb.Worksheets.Item(5).Name = "vSAT_TO-16502913000837517"

# Sheet 1 (GNG)
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16502912981622672.csv"
$ws1.Range("B3").Value = "GNG_stims-1650291298179743.csv"
$ws1.Range("B4").Value = "go_stims-1650291298181682.csv"
$ws1.Range("B5").Value = "GNG_stims-16502912981953251.csv"

# Sheet 2 (NB)
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-1650291299739449.csv"
$ws2.Range("B3").Value = "ZB-match_5-1650291298443221.csv"
$ws2.Range("B4").Value = "ZB-match_2-1650291298536497.csv"
$ws2.Range("B5").Value = "TB-16502912999071333.csv"
$ws2.Range("B6").Value = "OB-16502912985895226.csv"
$ws2.Range("B7").Value = "ZB-match_4-1650291298361807.csv"
$ws2.Range("B8").Value = "TB-16502912992538633.csv"
$ws2.Range("B9").Value = "OB-16502912987249017.csv"
$ws2.Range("B10").Value = "OB-16502912986166885.csv"

# Sheet 3 (RS)
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("B2").Value = "eyes closed"
$ws3.Range("B3").Value = "eyes open"

# Sheet 4 (TOL)
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16502912999484842.csv"
$ws4.Range("B3").Value = "ZM_stims-16502912999375386.csv"
$ws4.Range("B4").Value = "MM_stims-1650291299979706.csv"
$ws4.Range("B5").Value = "ZM_stims-1650291299949486.csv"
$ws4.Range("B6").Value = "MM_stims-16502912999949915.csv"
$ws4.Range("B7").Value = "ZM_stims-16502912999807098.csv"

# Sheet 5 (vSAT)
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "SAT_stims-16502913000252647.csv"
$ws5.Range("B3").Value = "SAT_stims-16502913000004067.csv"
$ws5.Range("B4").Value = "vSAT_stims-16502913000418155.csv"
$ws5.Range("B5").Value = "vSAT_stims-16502913000577142.csv"
